$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new column K values
$ws.Range("K3").Value = $null
$ws.Range("K4").Value = 2020
$ws.Range("K5").Value = 173

# Update selection to match the target active cell
$ws.Range("I18").Select()
